# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 ("R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 367
$wsOff.Range("C3").Value = 245
$wsOff.Range("D3").Value = 70
$wsOff.Range("E3").Value = 26
$wsOff.Range("F3").Value = 10

# DEF sheet - row 3 ("R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 458
$wsDef.Range("C3").Value = 328
$wsDef.Range("D3").Value = 110
$wsDef.Range("E3").Value = 57
$wsDef.Range("G3").Value = 8
